{"js": "const body = context.document.body;\n\nconst replacements = [\n  [\"31\u00d788=\", \"92\u00d780=\"],\n  [\"24\u00d775=\", \"19\u00d799=\"],\n  [\"91\u00d771=\", \"31\u00d777=\"],\n  [\"15\u00d787=\", \"45\u00d724=\"],\n  [\"96\u00d719=\", \"78\u00d731=\"],\n  [\"64\u00d758=\", \"13\u00d750=\"],\n  [\"53\u00d791=\", \"83\u00d797=\"],\n  [\"58\u00d731=\", \"86\u00d779=\"],\n  [\"81\u00d796=\", \"88\u00d762=\"],\n  [\"54\u00d753=\", \"64\u00d755=\"],\n  [\"65\u00d795=\", \"79\u00d724=\"],\n  [\"51\u00d761=\", \"24\u00d794=\"],\n  [\"55\u00d737=\", \"89\u00d789=\"],\n  [\"96\u00d790=\", \"20\u00d747=\"],\n  [\"85\u00d751=\", \"18\u00d771=\"],\n  [\"43\u00d745=\", \"13\u00d778=\"],\n  [\"38\u00d757=\", \"98\u00d714=\"],\n  [\"50\u00d785=\", \"64\u00d737=\"],\n  [\"92\u00d754=\", \"36\u00d759=\"],\n  [\"54\u00d775=\", \"34\u00d781=\"],\n  [\"41\u00d746=\", \"67\u00d743=\"],\n  [\"51\u00d766=\", \"45\u00d770=\"],\n  [\"34\u00d785=\", \"91\u00d737=\"],\n  [\"16\u00d743=\", \"46\u00d788=\"],\n  [\"40\u00d712=\", \"31\u00d711=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"31\u00d788=\", \"92\u00d780=\"),\n    @(\"24\u00d775=\", \"19\u00d799=\"),\n    @(\"91\u00d771=\", \"31\u00d777=\"),\n    @(\"15\u00d787=\", \"45\u00d724=\"),\n    @(\"96\u00d719=\", \"78\u00d731=\"),\n    @(\"64\u00d758=\", \"13\u00d750=\"),\n    @(\"53\u00d791=\", \"83\u00d797=\"),\n    @(\"58\u00d731=\", \"86\u00d779=\"),\n    @(\"81\u00d796=\", \"88\u00d762=\"),\n    @(\"54\u00d753=\", \"64\u00d755=\"),\n    @(\"65\u00d795=\", \"79\u00d724=\"),\n    @(\"51\u00d761=\", \"24\u00d794=\"),\n    @(\"55\u00d737=\", \"89\u00d789=\"),\n    @(\"96\u00d790=\", \"20\u00d747=\"),\n    @(\"85\u00d751=\", \"18\u00d771=\"),\n    @(\"43\u00d745=\", \"13\u00d778=\"),\n    @(\"38\u00d757=\", \"98\u00d714=\"),\n    @(\"50\u00d785=\", \"64\u00d737=\"),\n    @(\"92\u00d754=\", \"36\u00d759=\"),\n    @(\"54\u00d775=\", \"34\u00d781=\"),\n    @(\"41\u00d746=\", \"67\u00d743=\"),\n    @(\"51\u00d766=\", \"45\u00d770=\"),\n    @(\"34\u00d785=\", \"91\u00d737=\"),\n    @(\"16\u00d743=\", \"46\u00d788=\"),\n    @(\"40\u00d712=\", \"31\u00d711=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($null, $true, $false, $false, $null, $null, $true, 1, $null, $newText, 2) | Out-Null\n}\n"}
